$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "60.787.72"
Set-TextValue $ws.Range("E2") "  -0.11%  "
Set-TextValue $ws.Range("D3") "2.910.07"
Set-TextValue $ws.Range("E3") "  -0.01%  "
Set-TextValue $ws.Range("D5") "591.92"
Set-TextValue $ws.Range("E5") "  +1.05%  "
Set-TextValue $ws.Range("D6") "145.24"
Set-TextValue $ws.Range("E6") "  -0.72%  "
Set-TextValue $ws.Range("E7") "  -0.03%  "
Set-TextValue $ws.Range("E8") "  +0.70%  "
Set-TextValue $ws.Range("D9") "6.89"
Set-TextValue $ws.Range("E9") "  +1.63%  "
Set-TextValue $ws.Range("E10") "  -0.77%  "
Set-TextValue $ws.Range("D11") "0.438"
Set-TextValue $ws.Range("E11") "  -2.16%  "
Set-TextValue $ws.Range("E12") "  +0.26%  "
Set-TextValue $ws.Range("E13") "  -0.75%  "
Set-TextValue $ws.Range("E14") "  -0.84%  "
Set-TextValue $ws.Range("D15") "3.391.22"
Set-TextValue $ws.Range("E15") "  -0.07%  "
Set-TextValue $ws.Range("D16") "60.794.06"
Set-TextValue $ws.Range("E16") "  -0.08%  "
Set-TextValue $ws.Range("B17") "WrappedEther"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "2.912.03"
Set-TextValue $ws.Range("E17") "  +0.00%  "
Set-TextValue $ws.Range("B18") "Polkadot"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D18") "6.66"
Set-TextValue $ws.Range("E18") "  -1.37%  "
Set-TextValue $ws.Range("D19") "429.32"
Set-TextValue $ws.Range("E19") "  +0.51%  "
Set-TextValue $ws.Range("D20") "13.29"
Set-TextValue $ws.Range("E20") "  -2.02%  "
Set-TextValue $ws.Range("D21") "0.675"
Set-TextValue $ws.Range("E21") "  +0.49%  "
Set-TextValue $ws.Range("E22") "  -1.12%  "
Set-TextValue $ws.Range("D23") "81.34"
Set-TextValue $ws.Range("E23") "  +1.37%  "
Set-TextValue $ws.Range("E24") "  -0.28%  "
Set-TextValue $ws.Range("E25") "  -0.90%  "
Set-TextValue $ws.Range("E26") "  -0.03%  "
Set-TextValue $ws.Range("E28") "  +5.28%  "
Set-TextValue $ws.Range("E29") "  +0.00%  "
Set-TextValue $ws.Range("E30") "  -0.47%  "
Set-TextValue $ws.Range("D31") "7.01"
Set-TextValue $ws.Range("E31") "  -3.42%  "
Set-TextValue $ws.Range("D32") "26.45"
Set-TextValue $ws.Range("E32") "  +0.11%  "
Set-TextValue $ws.Range("E33") "  +0.78%  "
Set-TextValue $ws.Range("D34") "0.0₃0849"
Set-TextValue $ws.Range("E34") "  +0.94%  "
Set-TextValue $ws.Range("E35") "  +0.21%  "
Set-TextValue $ws.Range("D36") "5.60"
Set-TextValue $ws.Range("E36") "  -0.32%  "
Set-TextValue $ws.Range("E37") "  +0.64%  "
Set-TextValue $ws.Range("E38") "  -1.60%  "
Set-TextValue $ws.Range("E39") "  -1.57%  "
Set-TextValue $ws.Range("E40") "  -1.72%  "
Set-TextValue $ws.Range("E41") "  -2.39%  "
Set-TextValue $ws.Range("D42") "39.86"
Set-TextValue $ws.Range("E42") "  -4.29%  "
Set-TextValue $ws.Range("D43") "374.22"
Set-TextValue $ws.Range("E43") "  -0.74%  "
Set-TextValue $ws.Range("E44") "  -0.89%  "
Set-TextValue $ws.Range("D45") "2.698.37"
Set-TextValue $ws.Range("E45") "  +0.89%  "
Set-TextValue $ws.Range("D46") "132.03"
Set-TextValue $ws.Range("E46") "  -0.58%  "
Set-TextValue $ws.Range("D48") "23.67"
Set-TextValue $ws.Range("E48") "  -4.77%  "
Set-TextValue $ws.Range("E49") "  -0.49%  "
Set-TextValue $ws.Range("E50") "  -3.66%  "
Set-TextValue $ws.Range("E51") "  +0.73%  "
